$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.129.64"
$ws.Range("E2").Value = "  -0.14%  "
$ws.Range("D3").Value = "2.563.96"
$ws.Range("E3").Value = "  +0.64%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'584.37"
$ws.Range("E5").Value = "  +2.68%  "
$ws.Range("D6").Value = "'147.64"
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'0.597"
$ws.Range("E8").Value = "  +1.64%  "
$ws.Range("E9").Value = "  +2.91%  "
$ws.Range("E10").Value = "  +0.76%  "
$ws.Range("E11").Value = "  -0.09%  "
$ws.Range("D12").Value = "'0.357"
$ws.Range("E12").Value = "  +1.12%  "
$ws.Range("D13").Value = "'27.44"
$ws.Range("E13").Value = "  -0.10%  "
$ws.Range("D14").Value = "3.023.01"
$ws.Range("E14").Value = "  +0.63%  "
$ws.Range("D15").Value = "63.060.47"
$ws.Range("E15").Value = "  -0.11%  "
$ws.Range("E16").Value = "  +1.92%  "
$ws.Range("D17").Value = "2.554.12"
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("D18").Value = "'11.36"
$ws.Range("E18").Value = "  -1.04%  "
$ws.Range("D19").Value = "'343.46"
$ws.Range("E19").Value = "  +1.88%  "
$ws.Range("D20").Value = "'4.41"
$ws.Range("E20").Value = "  +2.59%  "
$ws.Range("E21").Value = "  +1.83%  "
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("E23").Value = "  -4.01%  "
$ws.Range("D24").Value = "'66.77"
$ws.Range("E24").Value = "  +1.97%  "
$ws.Range("D25").Value = "2.696.33"
$ws.Range("E25").Value = "  +0.74%  "
$ws.Range("E26").Value = "  +0.54%  "
$ws.Range("E27").Value = "  -0.78%  "
$ws.Range("D28").Value = "'8.08"
$ws.Range("E28").Value = "  +10.14%  "
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("D30").Value = "'8.43"
$ws.Range("E30").Value = "  -0.31%  "
$ws.Range("E31").Value = "  -2.43%  "
$ws.Range("E32").Value = "  +7.27%  "
$ws.Range("D33").Value = "0.0₃0825"
$ws.Range("E33").Value = "  +0.75%  "
$ws.Range("D34").Value = "'460.48"
$ws.Range("E34").Value = "  +11.88%  "
$ws.Range("D35").Value = "'1.62"
$ws.Range("E35").Value = "  +3.02%  "
$ws.Range("D36").Value = "'175.74"
$ws.Range("E36").Value = "  -1.39%  "
$ws.Range("D37").Value = "'0.407"
$ws.Range("E37").Value = "  +1.98%  "
$ws.Range("D38").Value = "'19.20"
$ws.Range("E38").Value = "  +0.97%  "
$ws.Range("D39").Value = "'4.54"
$ws.Range("E39").Value = "  +3.39%  "
$ws.Range("D41").Value = "'1.76"
$ws.Range("E41").Value = "  -1.03%  "
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("D43").Value = "'151.15"
$ws.Range("E43").Value = "  -1.60%  "
$ws.Range("E44").Value = "  +1.69%  "
$ws.Range("D45").Value = "'20.90"
$ws.Range("E45").Value = "  -0.24%  "
$ws.Range("D46").Value = "'0.0547"
$ws.Range("E46").Value = "  +5.10%  "
$ws.Range("E47").Value = "  +0.99%  "
$ws.Range("D48").Value = "'0.0975"
$ws.Range("E48").Value = "  +1.41%  "
$ws.Range("E49").Value = "  +0.53%  "
$ws.Range("E50").Value = "  -2.66%  "
$ws.Range("E51").Value = "  +0.23%  "
